$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.663.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.96"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5367"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3198"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07062"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.07"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7758"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07825"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.855.64"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.055"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.684.70"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.093.90"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.645"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.049"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.411"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.87"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.216"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.697"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.16"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.75"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.294"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08758"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.115"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04882"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7371"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.40%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.899"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.111"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.347"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01750"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4841"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9098"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.55"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.924"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.742"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4205"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.130"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1252"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05838"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8986"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.65%  "
